# Insert 3 new weekly rows of "Lechuga" price data into the subconjunto sheet.
# This pushes the existing rows 1079:1120 down to 1082:1123 (dimension grows
# from A1:R1120 to A1:R1123) and fills the newly opened rows 1079:1081 with
# the new week's records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the block (rows shift down).
$ws.Rows("1079:1081").Insert()

# New row 1079: Lechuga Conconina(o) Primera
$row1079 = New-Object 'object[,]' 1,18
$row1079[0,0]  = 7
$row1079[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row1079[0,2]  = "Ñuble"
$row1079[0,3]  = 45075
$row1079[0,4]  = 16
$row1079[0,5]  = 100112033
$row1079[0,6]  = "Lechuga"
$row1079[0,7]  = "Conconina(o)"
$row1079[0,8]  = "Primera"
$row1079[0,9]  = 100
$row1079[0,10] = 7000
$row1079[0,11] = 7000
$row1079[0,12] = 7000
$row1079[0,13] = "`$/caja 10 unidades"
$row1079[0,14] = "Región del Maule"
$row1079[0,15] = 700
$row1079[0,16] = 10
$row1079[0,17] = "Hortaliza"
$ws.Range("A1079:R1079").Value = $row1079

# New row 1080: Lechuga Conconina(o) Segunda
$row1080 = New-Object 'object[,]' 1,18
$row1080[0,0]  = 7
$row1080[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row1080[0,2]  = "Ñuble"
$row1080[0,3]  = 45075
$row1080[0,4]  = 16
$row1080[0,5]  = 100112033
$row1080[0,6]  = "Lechuga"
$row1080[0,7]  = "Conconina(o)"
$row1080[0,8]  = "Segunda"
$row1080[0,9]  = 80
$row1080[0,10] = 5000
$row1080[0,11] = 5000
$row1080[0,12] = 5000
$row1080[0,13] = "`$/caja 12 unidades"
$row1080[0,14] = "Región del Maule"
$row1080[0,15] = 417
$row1080[0,16] = 12
$row1080[0,17] = "Hortaliza"
$ws.Range("A1080:R1080").Value = $row1080

# New row 1081: Lechuga Escarola Primera
$row1081 = New-Object 'object[,]' 1,18
$row1081[0,0]  = 7
$row1081[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row1081[0,2]  = "Ñuble"
$row1081[0,3]  = 45075
$row1081[0,4]  = 16
$row1081[0,5]  = 100112033
$row1081[0,6]  = "Lechuga"
$row1081[0,7]  = "Escarola"
$row1081[0,8]  = "Primera"
$row1081[0,9]  = 100
$row1081[0,10] = 7000
$row1081[0,11] = 7000
$row1081[0,12] = 7000
$row1081[0,13] = "`$/caja 15 unidades"
$row1081[0,14] = "Región del Maule"
$row1081[0,15] = 467
$row1081[0,16] = 15
$row1081[0,17] = "Hortaliza"
$ws.Range("A1081:R1081").Value = $row1081
